$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top of the sheet, shifting all existing rows down by one
$ws.Rows.Item(1).Insert()

# Set the new header cell value
$ws.Range("A1").Value = "WordKey"
